$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new "picture" file-name column (E) for most of the rows (a set of
# image files uploaded alongside the workbook).
$ws.Range("E3").Value = "images/Piano.jpg"
$ws.Range("E4").Value = "images/A.jpg"
$ws.Range("E5").Value = "images/Alat_tulis.jpg"
$ws.Range("E6").Value = "images/Gen_Alpha.jpg"
$ws.Range("E7").Value = "images/Spiderman.jpg"
$ws.Range("E8").Value = "images/Bach.jpg"
$ws.Range("E9").Value = "images/Beethoven.jpg"
$ws.Range("E10").Value = "images/Liszt.jpg"
$ws.Range("E11").Value = "images/Chop.jpg"
$ws.Range("E12").Value = "images/Joshua.jpg"

# Widen column E so the new picture paths are fully visible.
$ws.Columns.Item(5).ColumnWidth = 51.5

# Match the author's final selection when the workbook was saved.
$ws.Range("E17").Select()
